# Rename the header labels to their uppercase / English equivalents.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "NAME"
$ws.Range("C1").Value = "TELEFONO"

# Re-write the CONCAT formulas without the stray leading space that the
# original file carried (" _xlfn.CONCAT(...)" -> "_xlfn.CONCAT(...)").
$ws.Range("A2").Formula = "=_xlfn.CONCAT(B2,C2)"
$ws.Range("A3").Formula = "=_xlfn.CONCAT(B3,C3)"
$ws.Range("A4").Formula = "=_xlfn.CONCAT(B4,C4)"
$ws.Range("A5").Formula = "=_xlfn.CONCAT(B5,C5)"

# Move the active selection to C6, matching where the author left off.
$null = $ws.Range("C6").Select()
